# Apply updated crypto price/volume figures to columns D (Price) and E (Volume(1h)).
# Source values come from the refreshed coinranking.com scrape (GitHub Actions run).
#
# Column D holds plain text (not numbers): some values look numeric ("1.002"), so we
# prefix them with a leading apostrophe - exactly like typing '1.002 into Excel - to
# keep Excel from reinterpreting them as the number 1.002 and dropping the formatting
# (trailing zeros, dot-as-thousands-separator, etc). Values that already contain two
# dots (e.g. "27.870.68") are never parsed as numbers by Excel, so no prefix is needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '27.870.68'
$ws.Range('E2').Value = '  +1.17%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.755.95'

# Row 4: TetherUSD
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.02%  '

# Row 5: BNB
$ws.Range('D5').Value = '''326.78'
$ws.Range('E5').Value = '  +0.70%  '

# Row 6: USDC
$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  -0.04%  '

# Row 7: XRP
$ws.Range('D7').Value = '''0.4583'
$ws.Range('E7').Value = '  +0.45%  '

# Row 8: Cardano
$ws.Range('D8').Value = '''0.3497'
$ws.Range('E8').Value = '  -1.56%  '

# Row 9: OKB
$ws.Range('E9').Value = '  +1.00%  '

# Row 10: Dogecoin
$ws.Range('D10').Value = '''0.07354'

# Row 11: Polygon
$ws.Range('E11').Value = '  -0.59%  '

# Row 12: BinanceUSD
$ws.Range('D12').Value = '''1.001'
$ws.Range('E12').Value = '  -0.01%  '

# Row 13: Solana
$ws.Range('D13').Value = '''20.54'
$ws.Range('E13').Value = '  -1.31%  '

# Row 14: Polkadot
$ws.Range('D14').Value = '''5.975'
$ws.Range('E14').Value = '  -0.48%  '

# Row 15: Chainlink
$ws.Range('D15').Value = '''7.144'
$ws.Range('E15').Value = '  -0.08%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '1.754.54'
$ws.Range('E16').Value = '  +0.19%  '

# Row 17: Litecoin
$ws.Range('D17').Value = '''91.74'
$ws.Range('E17').Value = '  -2.02%  '

# Row 18: ShibaInu
$ws.Range('E18').Value = '  -0.09%  '

# Row 19: TRON
$ws.Range('D19').Value = '''0.06417'
$ws.Range('E19').Value = '  +0.09%  '

# Row 20: Dai
$ws.Range('D20').Value = '''1.000'
$ws.Range('E20').Value = '  +0.01%  '

# Row 21: Avalanche
$ws.Range('D21').Value = '''16.82'
$ws.Range('E21').Value = '  -1.49%  '

# Row 22: Uniswap
$ws.Range('D22').Value = '''5.760'
$ws.Range('E22').Value = '  +0.40%  '

# Row 23: WrappedBTC
$ws.Range('D23').Value = '27.905.24'

# Row 24: Cosmos
$ws.Range('E24').Value = '  -0.95%  '

# Row 25: Toncoin
$ws.Range('D25').Value = '''2.159'
$ws.Range('E25').Value = '  +4.41%  '

# Row 26: Monero
$ws.Range('D26').Value = '''161.83'
$ws.Range('E26').Value = '  -2.42%  '

# Row 27: EthereumClassic
$ws.Range('D27').Value = '''20.03'
$ws.Range('E27').Value = '  -0.51%  '

# Row 28: WrappedliquidstakedEther2.0
$ws.Range('D28').Value = '1.956.97'
$ws.Range('E28').Value = '  +0.19%  '

# Row 29: LidoDAOToken
$ws.Range('D29').Value = '''2.141'
$ws.Range('E29').Value = '  +0.41%  '

# Row 30: BitcoinCash
$ws.Range('D30').Value = '''122.90'
$ws.Range('E30').Value = '  -2.04%  '

# Row 31: ImmutableX
$ws.Range('D31').Value = '''1.064'
$ws.Range('E31').Value = '  -1.44%  '

# Row 32: Stellar
$ws.Range('D32').Value = '''0.09246'
$ws.Range('E32').Value = '  +0.31%  '

# Row 33: HuobiToken
$ws.Range('D33').Value = '''3.668'
$ws.Range('E33').Value = '  +0.20%  '

# Row 34: Filecoin
$ws.Range('D34').Value = '''5.530'
$ws.Range('E34').Value = '  +0.07%  '

# Row 35: Aptos
$ws.Range('E35').Value = '  -0.56%  '

# Row 36: VeChain
$ws.Range('E36').Value = '  -0.76%  '

# Row 37: Hedera
$ws.Range('E37').Value = '  +1.14%  '

# Row 38: Algorand
$ws.Range('D38').Value = '''0.2061'
$ws.Range('E38').Value = '  -1.50%  '

# Row 39: InternetComputer(DFINITY)
$ws.Range('D39').Value = '''4.898'
$ws.Range('E39').Value = '  -0.48%  '

# Row 40: TheSandbox
$ws.Range('E40').Value = '  -1.76%  '

# Row 41: TrustWalletToken
$ws.Range('D41').Value = '''1.178'
$ws.Range('E41').Value = '  -0.35%  '

# Row 42: FraxShare
$ws.Range('D42').Value = '''7.771'
$ws.Range('E42').Value = '  -0.20%  '

# Row 43: WEMIXTOKEN
$ws.Range('D43').Value = '''1.344'
$ws.Range('E43').Value = '  -3.04%  '

# Row 44: EnergySwap
$ws.Range('E44').Value = '  -0.67%  '

# Row 45: PancakeSwap
$ws.Range('D45').Value = '''3.720'
$ws.Range('E45').Value = '  +0.05%  '

# Row 46: Decentraland
$ws.Range('D46').Value = '''0.5764'
$ws.Range('E46').Value = '  -1.72%  '

# Row 47: Quant
$ws.Range('D47').Value = '''123.19'
$ws.Range('E47').Value = '  +0.91%  '

# Row 48: NEARProtocol
$ws.Range('E48').Value = '  -0.69%  '

# Row 49: Cronos
$ws.Range('D49').Value = '''0.06801'
$ws.Range('E49').Value = '  -1.43%  '

# Row 50: EOS
$ws.Range('D50').Value = '''1.121'
$ws.Range('E50').Value = '  -0.76%  '

# Row 51: Aave
$ws.Range('D51').Value = '''72.01'
$ws.Range('E51').Value = '  -0.33%  '
